$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1 - copy formatting from the adjacent "sum"
# header (G1) so it matches the other header cells, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Row 2 value for the new Save column.
$ws.Range("H2").Value = 1
